# "Routing Master" sheet holds a single example/test-data row (row 2) used
# by the QA automation templates. This commit refreshes the sample
# "Item Number" (col B) and its Salesforce record Id (col D) to a newly
# generated engineering-item test record, as part of adding the
# "SO To inspection order" / "SO to RMA Receipt" RMA test cases.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Routing Master")

# B2: sample Item Number used by the RMA test cases
$ws.Range("B2").Value = "Pro-PEItem-LQRXE"

# D2: the corresponding Salesforce Id for that item
$ws.Range("D2").Value = "a345f000000u5GYAAY"

# Columns B and D are best-fit/auto-sized; nudge them back to fit the new
# (slightly different-width) text, same as Excel does automatically when
# the sheet is re-saved after the cell text changes.
$ws.Columns.Item(2).ColumnWidth = 16.833333333333336
$ws.Columns.Item(4).ColumnWidth = 20.333333333333336
